$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D ("Price") holds numeric-looking text (e.g. "66.069.53",
# "1.00", "0.0000330") that must stay as literal text -- Excel would
# otherwise silently coerce it to a Double and drop the original
# formatting/trailing zeros. Force the cell to Text format ("@") right
# before assigning each such value.

# --- Price (D) and Volume(1h) (E) updates for unchanged coin rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.069.53"
$ws.Range("E2").Value = "  +4.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.816.85"
$ws.Range("E3").Value = "  +8.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "429.64"
$ws.Range("E5").Value = "  +9.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.47"
$ws.Range("E6").Value = "  +11.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  +4.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.748"
$ws.Range("E9").Value = "  +10.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +7.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000330"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.43"
$ws.Range("E12").Value = "  +12.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.73"
$ws.Range("E13").Value = "  +17.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.398.27"
$ws.Range("E14").Value = "  +7.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.20"
$ws.Range("E15").Value = "  +18.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.774.58"
$ws.Range("E17").Value = "  +6.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.27"
$ws.Range("E18").Value = "  +8.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.13"
$ws.Range("E19").Value = "  +11.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.178.12"
$ws.Range("E20").Value = "  +4.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "416.49"
$ws.Range("E21").Value = "  +6.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.36"
$ws.Range("E22").Value = "  +10.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.51"
$ws.Range("E24").Value = "  +7.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "37.19"
$ws.Range("E25").Value = "  +10.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.13"
$ws.Range("E26").Value = "  +48.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.34"
$ws.Range("E27").Value = "  +11.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("E28").Value = "  +13.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.41"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "14.13"
$ws.Range("E30").Value = "  +19.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "702.93"
$ws.Range("E31").Value = "  +4.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0479"
$ws.Range("E39").Value = "  +9.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  +46.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.329"
$ws.Range("E45").Value = "  +19.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.42"
$ws.Range("E46").Value = "  +11.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.19"
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.95"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.85"
$ws.Range("E51").Value = "  +7.22%  "

# --- Volume(1h) (E) only updates ---
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("E23").Value = "  +14.38%  "
$ws.Range("E32").Value = "  +16.35%  "
$ws.Range("E33").Value = "  +5.32%  "
$ws.Range("E38").Value = "  +5.09%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  +7.38%  "

# --- Row reordering / swaps (coin pairs swapped position with updated values) ---
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "39.84"
$ws.Range("E34").Value = "  +8.11%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.154"
$ws.Range("E36").Value = "  +3.09%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.77"
$ws.Range("E37").Value = "  +41.53%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0691"
$ws.Range("E41").Value = "  +14.77%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("E42").Value = "  +8.13%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.08"
$ws.Range("E48").Value = "  +6.50%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.65"
$ws.Range("E50").Value = "  +6.48%  "
